$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2499.1482
$ws.Range("I15").Value = 2499.1482
$ws.Range("K15").Value = 7497.444600000001
$ws.Range("M15").Value = -7328.444600000001

$ws.Range("H43").Value = 373542.28
$ws.Range("I43").Value = 570
$ws.Range("J43").Value = 586669.3
$ws.Range("K43").Value = 570
$ws.Range("L43").Value = 586669.3
$ws.Range("M43").Value = -501
$ws.Range("N43").Value = -586807.3

$ws.Range("H87").Value = 52498.875
$ws.Range("J87").Value = 52498.875
$ws.Range("L87").Value = 52498.875
$ws.Range("N87").Value = -54994.875

$ws.Range("H90").Value = 52498.875
$ws.Range("J90").Value = 52498.875
$ws.Range("L90").Value = 157496.625
$ws.Range("N90").Value = -169976.625

$ws.Range("H92").Value = 1618
$ws.Range("I92").Value = 959
$ws.Range("J92").Value = 2496.6667
$ws.Range("K92").Value = 959
$ws.Range("L92").Value = 2496.6667
$ws.Range("M92").Value = 289
$ws.Range("N92").Value = -4992.6667

$ws.Range("H125").Value = 83334400
$ws.Range("J125").Value = 1280.2
$ws.Range("L125").Value = 11521.8
$ws.Range("N125").Value = -16441.8

$ws.Range("H132").Value = 1516.8823
$ws.Range("I132").Value = 1486.5
$ws.Range("K132").Value = 4459.5
$ws.Range("M132").Value = -1929.5

$ws.Range("H138").Value = 2523.9092
$ws.Range("I138").Value = 2582.6
$ws.Range("J138").Value = 2475
$ws.Range("K138").Value = 7747.799999999999
$ws.Range("L138").Value = 7425
$ws.Range("M138").Value = -2607.799999999999
$ws.Range("N138").Value = -17705

$ws.Range("H141").Value = 23812370
$ws.Range("I141").Value = 23812370
$ws.Range("K141").Value = 71437110
$ws.Range("M141").Value = -71431930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9488.799999999999
$ws.Range("J45").Value = 10860.75
$ws.Range("L45").Value = 10860.75
$ws.Range("N45").Value = -11614.75

$ws.Range("H61").Value = 3303.8367
$ws.Range("I61").Value = 2185.975
$ws.Range("K61").Value = 2185.975
$ws.Range("M61").Value = -1973.975

$ws.Range("H122").Value = 17863.643
$ws.Range("I122").Value = 19826.455
$ws.Range("K122").Value = 59479.36500000001
$ws.Range("M122").Value = -57029.36500000001

$ws.Range("H132").Value = 3063.775
$ws.Range("I132").Value = 1824.069
$ws.Range("K132").Value = 5472.207
$ws.Range("M132").Value = -2942.207

$ws.Range("H136").Value = 3303.8367
$ws.Range("I136").Value = 2185.975
$ws.Range("K136").Value = 6557.924999999999
$ws.Range("M136").Value = -4007.924999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1745.125
$ws.Range("I105").Value = 1316.2142
$ws.Range("J105").Value = 4747.5
$ws.Range("K105").Value = 1316.2142
$ws.Range("L105").Value = 4747.5
$ws.Range("M105").Value = 430.7858000000001
$ws.Range("N105").Value = -8241.5

$ws.Range("H134").Value = 5147.875
$ws.Range("I134").Value = 2330.3845
$ws.Range("K134").Value = 6991.1535
$ws.Range("M134").Value = -4456.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7414432.5
$ws.Range("I31").Value = 2577.913
$ws.Range("J31").Value = 15163190
$ws.Range("K31").Value = 2577.913
$ws.Range("L31").Value = 15163190
$ws.Range("M31").Value = -2282.913
$ws.Range("N31").Value = -15163780

$ws.Range("H34").Value = 7414432.5
$ws.Range("I34").Value = 2577.913
$ws.Range("J34").Value = 15163190
$ws.Range("K34").Value = 2577.913
$ws.Range("L34").Value = 15163190
$ws.Range("M34").Value = -2375.913
$ws.Range("N34").Value = -15163594

$ws.Range("H58").Value = 18526900
$ws.Range("I58").Value = 62502492
$ws.Range("J58").Value = 10860.737
$ws.Range("K58").Value = 62502492
$ws.Range("L58").Value = 10860.737
$ws.Range("M58").Value = -62502289
$ws.Range("N58").Value = -11266.737

$ws.Range("H103").Value = 44788.2
$ws.Range("I103").Value = 19999.5
$ws.Range("J103").Value = 61314
$ws.Range("K103").Value = 19999.5
$ws.Range("L103").Value = 61314
$ws.Range("M103").Value = -18827.5
$ws.Range("N103").Value = -63658

$ws.Range("H105").Value = 5105803.5
$ws.Range("I105").Value = 7937550
$ws.Range("J105").Value = 8658.799999999999
$ws.Range("K105").Value = 7937550
$ws.Range("L105").Value = 8658.799999999999
$ws.Range("M105").Value = -7935803
$ws.Range("N105").Value = -12152.8

$ws.Range("H132").Value = 10261812
$ws.Range("I132").Value = 2400.8262
$ws.Range("K132").Value = 7202.4786
$ws.Range("M132").Value = -4672.4786

$ws.Range("H134").Value = 7093.7144
$ws.Range("I134").Value = 1567
$ws.Range("J134").Value = 10669.823
$ws.Range("K134").Value = 4701
$ws.Range("L134").Value = 32009.469
$ws.Range("M134").Value = -2166
$ws.Range("N134").Value = -37079.469

$ws.Range("H136").Value = 18526900
$ws.Range("I136").Value = 62502492
$ws.Range("J136").Value = 10860.737
$ws.Range("K136").Value = 187507476
$ws.Range("L136").Value = 32582.211
$ws.Range("M136").Value = -187504926
$ws.Range("N136").Value = -37682.211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 168.25
$ws.Range("I26").Value = 34.8
$ws.Range("K26").Value = 104.4
$ws.Range("M26").Value = 183.6

$ws.Range("H131").Value = 958.087
$ws.Range("I131").Value = 742.2
$ws.Range("J131").Value = 2397.3333
$ws.Range("K131").Value = 2226.6
$ws.Range("L131").Value = 7191.999899999999
$ws.Range("M131").Value = 2813.4
$ws.Range("N131").Value = -17271.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8319.409
$ws.Range("I70").Value = 6821.2
$ws.Range("J70").Value = 9567.916999999999
$ws.Range("K70").Value = 6821.2
$ws.Range("L70").Value = 9567.916999999999
$ws.Range("M70").Value = -6551.2
$ws.Range("N70").Value = -10107.917

$ws.Range("H73").Value = 8319.409
$ws.Range("I73").Value = 6821.2
$ws.Range("J73").Value = 9567.916999999999
$ws.Range("K73").Value = 6821.2
$ws.Range("L73").Value = 9567.916999999999
$ws.Range("M73").Value = -5885.2
$ws.Range("N73").Value = -11439.917

$ws.Range("H80").Value = 4499.5
$ws.Range("I80").Value = 7000
$ws.Range("J80").Value = 3666
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 3666
$ws.Range("M80").Value = -6002
$ws.Range("N80").Value = -5662

$ws.Range("H83").Value = 4499.5
$ws.Range("I83").Value = 7000
$ws.Range("J83").Value = 3666
$ws.Range("K83").Value = 35000
$ws.Range("L83").Value = 18330
$ws.Range("M83").Value = -30008
$ws.Range("N83").Value = -28314

$ws.Range("H97").Value = 908.1429000000001
$ws.Range("I97").Value = 844.9
$ws.Range("J97").Value = 1287.6
$ws.Range("K97").Value = 844.9
$ws.Range("L97").Value = 1287.6
$ws.Range("M97").Value = -348.9
$ws.Range("N97").Value = -2279.6

$ws.Range("H132").Value = 1908.6774
$ws.Range("I132").Value = 1555.9584
$ws.Range("J132").Value = 3118
$ws.Range("K132").Value = 4667.8752
$ws.Range("L132").Value = 9354
$ws.Range("M132").Value = -2137.8752
$ws.Range("N132").Value = -14414

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 568.4286
$ws.Range("I16").Value = 329.83334
$ws.Range("K16").Value = 329.83334
$ws.Range("M16").Value = -159.83334

$ws.Range("H22").Value = 1145.4412
$ws.Range("I22").Value = 694.1923
$ws.Range("J22").Value = 2612
$ws.Range("K22").Value = 694.1923
$ws.Range("L22").Value = 2612
$ws.Range("M22").Value = -399.1923
$ws.Range("N22").Value = -3202

$ws.Range("H27").Value = 1145.4412
$ws.Range("I27").Value = 694.1923
$ws.Range("J27").Value = 2612
$ws.Range("K27").Value = 694.1923
$ws.Range("L27").Value = 2612
$ws.Range("M27").Value = -587.1923
$ws.Range("N27").Value = -2826

$ws.Range("H46").Value = 2655495.8
$ws.Range("I46").Value = 17241768
$ws.Range("J46").Value = 3446.0908
$ws.Range("K46").Value = 17241768
$ws.Range("L46").Value = 3446.0908
$ws.Range("M46").Value = -17241580
$ws.Range("N46").Value = -3822.0908

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H132").Value = 11116120
$ws.Range("I132").Value = 18521534
$ws.Range("J132").Value = 7999.9443
$ws.Range("K132").Value = 55564602
$ws.Range("L132").Value = 23999.8329
$ws.Range("M132").Value = -55562072
$ws.Range("N132").Value = -29059.8329

$ws.Range("H136").Value = 6671.3184
$ws.Range("I136").Value = 2461.5
$ws.Range("K136").Value = 7384.5
$ws.Range("M136").Value = -4834.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 15469
$ws.Range("J50").Value = 15469
$ws.Range("L50").Value = 15469
$ws.Range("N50").Value = -16731

$ws.Range("H54").Value = 16875

$ws.Range("H107").Value = 12346296
$ws.Range("I107").Value = 474
$ws.Range("J107").Value = 30303856
$ws.Range("K107").Value = 1422
$ws.Range("L107").Value = 90911568
$ws.Range("M107").Value = 498
$ws.Range("N107").Value = -90915408

$ws.Range("H122").Value = 105721.8
$ws.Range("I122").Value = 139790.2
$ws.Range("J122").Value = 6923.4
$ws.Range("K122").Value = 419370.6
$ws.Range("L122").Value = 20770.2
$ws.Range("M122").Value = -416920.6
$ws.Range("N122").Value = -25670.2

$ws.Range("H132").Value = 16680495
$ws.Range("I132").Value = 29419092
$ws.Range("K132").Value = 88257276
$ws.Range("M132").Value = -88254746

$ws.Range("H136").Value = 23283892
$ws.Range("I136").Value = 43479610
$ws.Range("K136").Value = 130438830
$ws.Range("M136").Value = -130436280
